$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title
#    (a bold "Meta description" run followed by a plain run with the
#    description text). It was paragraph #2 in the original document.
$metaPara = $d.Paragraphs.Item(2)
if (-not $metaPara.Range.Text.Contains("Meta description")) {
    throw "Expected paragraph 2 to contain 'Meta description', found: $($metaPara.Range.Text)"
}
$metaPara.Range.Delete()

# 2. At the end of the document, the closing paragraph held the italic
#    "Create a feature image..." prompt. It gets replaced by two
#    paragraphs: a new bold title paragraph, followed by the same
#    (italic) paragraph but now holding the meta-description text
#    instead of the image prompt.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
if (-not $lastPara.Range.Text.Contains("Create a feature image")) {
    throw "Expected last paragraph to contain the image prompt, found: $($lastPara.Range.Text)"
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 40 Joker Staxx: 40 Lines for Free | Game Review</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the simple gameplay and high volatility of 40 Joker Staxx: 40 Lines. Play for free and find out if this classic slot game is right for you.</w:t></w:r></w:p>'

$lastPara.Range.InsertXML($xml)
